$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to store as text,
# matching the original inline-string cell type in the workbook.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '66.121.55'
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').Value = '3.273.59'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue 'D5' '575.52'
$ws.Range('E5').Value = '  -0.47%  '
Set-TextValue 'D6' '179.15'
$ws.Range('E6').Value = '  -3.96%  '
Set-TextValue 'D7' '0.625'
$ws.Range('E7').Value = '  +3.01%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('E11').Value = '  -1.69%  '
$ws.Range('D12').Value = '3.837.82'
$ws.Range('E12').Value = '  -1.44%  '
$ws.Range('E13').Value = '  -3.81%  '
$ws.Range('D14').Value = '66.176.62'
$ws.Range('E14').Value = '  -2.04%  '
Set-TextValue 'D15' '26.51'
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').Value = '3.267.53'
$ws.Range('E17').Value = '  -1.85%  '
Set-TextValue 'D18' '434.03'
$ws.Range('E18').Value = '  -2.35%  '
Set-TextValue 'D19' '5.57'
$ws.Range('E19').Value = '  -2.09%  '
Set-TextValue 'D20' '13.16'
$ws.Range('E20').Value = '  -3.07%  '
Set-TextValue 'D21' '7.41'
$ws.Range('E21').Value = '  -4.21%  '
Set-TextValue 'D22' '71.97'
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '3.410.88'
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('E25').Value = '  -2.00%  '
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('E27').Value = '  -5.39%  '
Set-TextValue 'D28' '8.91'
$ws.Range('E28').Value = '  -1.76%  '
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('E30').Value = '  -2.06%  '
Set-TextValue 'D31' '22.30'
$ws.Range('E31').Value = '  -2.81%  '
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('E33').Value = '  -3.68%  '
Set-TextValue 'D34' '6.60'
$ws.Range('E34').Value = '  -3.70%  '
$ws.Range('E35').Value = '  -4.78%  '
Set-TextValue 'D36' '156.86'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('E37').Value = '  -5.97%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D38' '26.51'
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D39' '1.79'
$ws.Range('E39').Value = '  -3.38%  '
$ws.Range('D40').Value = '2.775.50'
$ws.Range('E40').Value = '  -0.36%  '
Set-TextValue 'D41' '0.775'
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('E42').Value = '  -3.93%  '
Set-TextValue 'D43' '40.24'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  -3.68%  '
Set-TextValue 'D45' '0.0656'
$ws.Range('E45').Value = '  -2.49%  '
Set-TextValue 'D46' '2.30'
$ws.Range('E46').Value = '  -4.17%  '
Set-TextValue 'D47' '320.69'
$ws.Range('E47').Value = '  -1.87%  '
Set-TextValue 'D48' '23.31'
$ws.Range('E48').Value = '  -5.98%  '
Set-TextValue 'D49' '0.0267'
$ws.Range('E49').Value = '  -2.38%  '
Set-TextValue 'D50' '0.102'
$ws.Range('E50').Value = '  +1.92%  '
Set-TextValue 'D51' '1.00'
$ws.Range('E51').Value = '  +0.02%  '
